# The author re-labeled the "Rendimento médio na extração vegetal" variable
# (column B) to "Valor médio unitário na extração vegetal" across every row
# of the sheet that used it (220 rows in total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

$ws.Cells.Replace(
    "Rendimento médio na extração vegetal",
    "Valor médio unitário na extração vegetal",
    $xlWhole
)
